$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for rows 2-6 (columns D through AJ)
$cols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ")
$vals = @(442,-36,-36,-46,-46,-46,0,420,219,201,202,-1,41,-20,-18,56,4,-24,93,-8.17,-10.45,-20.81,-11.36,109,416.54,-576,-9.68,2491,2.24,0,0,0,8147358)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value = $vals[$i]
}

$cols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ")
$vals = @(604,45,45,21,21,21,0,464,178,287,287,-1,47,36,-39,13,8,28,44,7.49,3.45,8.460000000000001,4.71,62.08,527.3099999999999,233,28.81,3059,2.19,0,0,0,9427358)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "3").Value = $vals[$i]
}

$cols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ")
$vals = @(715,48,48,45,35,35,0,643,279,364,364,0,50,78,-158,99,31,46,103,6.76,4.85,10.6,6.27,76.61,632.86,362,47.6,3640,4.74,0,0,0,10055292)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "4").Value = $vals[$i]
}

$cols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ")
$vals = @(799,60,60,61,44,43,1,685,174,511,510,1,56,18,-20,-1,15,3,2,7.52,5.55,9.880000000000001,6.68,34.13,799.59,393,29.12,4533,2.53,0,0,0,11288324)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "5").Value = $vals[$i]
}

$cols = @("D","E","F","G","H","I","K","L","M","N","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ")
$vals = @(891,50,50,39,45,45,982,413,568,567,58,42,-93,223,9,33,193,5.6,5.05,8.289999999999999,5.39,72.73,895.51,389,19.04,4922,1.5,0,0,0,11570324)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "6").Value = $vals[$i]
}

# Rows 7-9: clear all data columns (D:AJ), keep A, B, C
$ws.Range("D7:AJ9").ClearContents()

Write-Host "Edit complete"